$d = $word.ActiveDocument
$d.Content.Find.Execute("Hành vi", $true, $false, $false, $false, $false, $true, 1, $false, "function", 2)
$d.Content.Find.Execute("Thuộc tính", $true, $false, $false, $false, $false, $true, 1, $false, "property", 2)
